$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D for the newest fiscal year (2018-12-31), shifting
# all existing year columns (D:K) one column to the right (E:L), and copy the
# number formatting from the (old) column D -- now column E -- into the new column D.
$ws.Columns("D").Insert(-4161)
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the newly inserted column D with the FY2018 figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 824800
$ws.Range("D9").Value2 = 682700
$ws.Range("D10").Value2 = 142100
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 8900
$ws.Range("D15").Value2 = 388900
$ws.Range("D17").Value2 = 1112600
$ws.Range("D18").Value2 = -287800
$ws.Range("D20").Value2 = 45100
$ws.Range("D21").Value2 = 146200
$ws.Range("D22").Value2 = 156300
$ws.Range("D23").Value2 = -399000
$ws.Range("D24").Value2 = -44200
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = -354800
$ws.Range("D27").Value2 = -354800
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 7400
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -45100
$ws.Range("D33").Value2 = -347400
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = -347400
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 1026700
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 259900
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 13700
$ws.Range("D46").Value2 = 1300300
$ws.Range("D47").Value2 = 497600
$ws.Range("D48").Value2 = 6201000
$ws.Range("D49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 118800
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 8117700
$ws.Range("D57").Value2 = 122300
$ws.Range("D58").Value2 = 201200
$ws.Range("D59").Value2 = 130100
$ws.Range("D60").Value2 = 453600
$ws.Range("D61").Value2 = 2309700
$ws.Range("D62").Value2 = 319400
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 3082700
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 3810500
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 5035000
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = -347400
$ws.Range("D83").Value2 = 388900
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = -160100
$ws.Range("D91").Value2 = -240000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -134100
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -11200
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = -305400

# A couple of cells were corrected (not just shifted) as part of this update.
$ws.Range("F89").Value2 = 929600
$ws.Range("F100").Value2 = -46900
